# Weekly update: insert a new week's worth of Alcachofa price records
# (Comercializadora del Agro de Limarí) at the top of the data block and
# push the rest of the historical rows down by 3 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at row 127 (pushes existing rows 127.. down to 130..)
$ws.Range("A127:A129").EntireRow.Insert()

# New week's data (fecha serial 44468 == 2021-09-29)
$ws.Cells.Item(127, 1).Value  = 2
$ws.Cells.Item(127, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(127, 3).Value  = "Coquimbo"
$ws.Cells.Item(127, 4).Value  = 44468
$ws.Cells.Item(127, 5).Value  = 4
$ws.Cells.Item(127, 6).Value  = 100112013
$ws.Cells.Item(127, 7).Value  = "Alcachofa"
$ws.Cells.Item(127, 8).Value  = "Argentina(o)"
$ws.Cells.Item(127, 9).Value  = "Primera"
$ws.Cells.Item(127, 10).Value = 1000
$ws.Cells.Item(127, 11).Value = 7000
$ws.Cells.Item(127, 12).Value = 8000
$ws.Cells.Item(127, 13).Value = 7500
$ws.Cells.Item(127, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(127, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(127, 16).Value = 150
$ws.Cells.Item(127, 17).Value = 50
$ws.Cells.Item(127, 18).Value = "Hortaliza"

$ws.Cells.Item(128, 1).Value  = 2
$ws.Cells.Item(128, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(128, 3).Value  = "Coquimbo"
$ws.Cells.Item(128, 4).Value  = 44468
$ws.Cells.Item(128, 5).Value  = 4
$ws.Cells.Item(128, 6).Value  = 100112013
$ws.Cells.Item(128, 7).Value  = "Alcachofa"
$ws.Cells.Item(128, 8).Value  = "Argentina(o)"
$ws.Cells.Item(128, 9).Value  = "Segunda"
$ws.Cells.Item(128, 10).Value = 1100
$ws.Cells.Item(128, 11).Value = 5000
$ws.Cells.Item(128, 12).Value = 6000
$ws.Cells.Item(128, 13).Value = 5500
$ws.Cells.Item(128, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(128, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(128, 16).Value = 79
$ws.Cells.Item(128, 17).Value = 70
$ws.Cells.Item(128, 18).Value = "Hortaliza"

$ws.Cells.Item(129, 1).Value  = 2
$ws.Cells.Item(129, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(129, 3).Value  = "Coquimbo"
$ws.Cells.Item(129, 4).Value  = 44468
$ws.Cells.Item(129, 5).Value  = 4
$ws.Cells.Item(129, 6).Value  = 100112013
$ws.Cells.Item(129, 7).Value  = "Alcachofa"
$ws.Cells.Item(129, 8).Value  = "Española"
$ws.Cells.Item(129, 9).Value  = "Primera"
$ws.Cells.Item(129, 10).Value = 1300
$ws.Cells.Item(129, 11).Value = 9000
$ws.Cells.Item(129, 12).Value = 10000
$ws.Cells.Item(129, 13).Value = 9500
$ws.Cells.Item(129, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(129, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(129, 16).Value = 317
$ws.Cells.Item(129, 17).Value = 30
$ws.Cells.Item(129, 18).Value = "Hortaliza"
